$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" - a new handoff xliff (b.*.xlf) was produced
# for b.md, so the localization-status report's "b.md" rows move from
# "Handed back: in sync with en-US" to "Ready for handoff", pick up the new
# handoff file names/timestamps, and surface a version-mismatch error.
# ---------------------------------------------------------------------------

$statusNew      = "Ready for handoff"
$hoDateOverview = "2016-08-28 00:35:39"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5237317f0d835b0a9dd99677f083d0a40788af52/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4034ea7e90240ce8c997c498e78910e48214e5d1/e2e/b.md."

# ----------------------------- Overview sheet ------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = $hoDateOverview

# -------------------------------- zh-cn sheet -------------------------------
# NOTE: a leading apostrophe forces "False" to be stored as literal text
# (matching the source report's column type) instead of auto-converting to
# the Excel Boolean FALSE the way a bare "False"/"True" entry normally would.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusNew
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-28 00:35:34"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# -------------------------------- de-de sheet -------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-28 00:35:39"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
